# Edit 1: Merge "Celaldoğan" + " Güneş | 05200000067" runs into a single run,
# dropping the now-stray proofErr spell-check markers.
$d = $word.ActiveDocument
$d.Content.Find.Execute("Celaldoğan Güneş | 05200000067", $true, $false, $false, $false, $false, $true, 1, $false, "Celaldoğan Güneş | 05200000067", 2)

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Celaldoğan*") {
        $xml = $p.Range.WordOpenXML
        $p.Range.InsertXML($xml)
    }
}

# Edit 2: Split the "Sırasıyla aracın ..." run so a new clause about the driver
# riding a motorcycle/passenger vehicle is inserted before the disability clause.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Sırasıyla aracın*") {
        $frag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="32E99BE6" w14:textId="7C1AC3B8" w:rsidR="00664608" w:rsidRPr="00664608" w:rsidRDefault="00BD3A7F" w:rsidP="00BD3A7F"><w:pPr><w:pStyle w:val="ListeParagraf"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:r w:rsidR="00E93148" w:rsidRPr="00C5588A"><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Sırasıyla aracın plakasını, sınıf kodunu, aracın ağırlığını (kilogram), aracın otoparkta kaldığı süreyi (dakika), sürücünün adını soyadını</w:t></w:r><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> ve sürücü </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>motorsiklet</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> veya binek araç sürüyorsa</w:t></w:r><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> sürücünün engellilik / gazilik özel durumunu “E,</w:t></w:r><w:r w:rsidR="00C5588A"><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00E93148" w:rsidRPr="00C5588A"><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>e,</w:t></w:r><w:r w:rsidR="00C5588A"><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00E93148" w:rsidRPr="00C5588A"><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>G,</w:t></w:r><w:r w:rsidR="00C5588A"><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00E93148" w:rsidRPr="00C5588A"><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>g” harflerini kullanarak girin.</w:t></w:r><w:r w:rsidR="005D1EEC" w:rsidRPr="00C5588A"><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> Eğer özel durumu yoksa “y” veya “Y” harflerinden birini girin.</w:t></w:r><w:r w:rsidR="005D1EEC" w:rsidRPr="00C5588A"><w:rPr><w:noProof/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($frag)
    }
}

# Edit 3: Mark the picture run (the one holding the 6645910x1341755 drawing) as NoProof.
$shp = $d.InlineShapes.Item(7)
$shp.Range.NoProofing = 1
